$wb = $excel.ActiveWorkbook

# --- "|Toets 1" sheet: add a new "aantal vragen:" / count column ---
$toets1 = $wb.Worksheets.Item("|Toets 1")
$toets1.Range("G1").Value = "aantal vragen:"
$toets1.Range("H1").Value = 6

# --- "|Toets 2" sheet: add a new "Aantal vragen:" / count column ---
$toets2 = $wb.Worksheets.Item("|Toets 2")
$toets2.Range("G1").Value = "Aantal vragen:"
$toets2.Range("H1").Value = 6

# Data import finished up on the "|Toets 2" sheet: make it the active tab,
# with the selection resting on K14.
$toets2.Activate() | Out-Null
$toets2.Range("K14").Select() | Out-Null
